$wb = $excel.ActiveWorkbook

# --- Citywide Totals ---
$ws = $wb.Worksheets.Item("Citywide Totals")
$ws.Range("F2").Value = 57   # was 56
$ws.Range("B3").Value = 54   # was 53
$ws.Range("G3").Value = 79   # was 77
$ws.Range("J3").Value = 134   # was 131
$ws.Range("B6").Value = 241   # was 240
$ws.Range("C6").Value = 309   # was 307
$ws.Range("E6").Value = 272   # was 271
$ws.Range("F6").Value = 361   # was 360
$ws.Range("G6").Value = 314   # was 312
$ws.Range("H6").Value = 281   # was 278
$ws.Range("K6").Value = 332   # was 331
$ws.Range("B7").Value = 329   # was 327
$ws.Range("C7").Value = 413   # was 411
$ws.Range("E7").Value = 420   # was 419
$ws.Range("F7").Value = 511   # was 509
$ws.Range("G7").Value = 459   # was 455
$ws.Range("H7").Value = 441   # was 438
$ws.Range("J7").Value = 478   # was 475
$ws.Range("K7").Value = 588   # was 587

# --- Garfield Park ---
$ws = $wb.Worksheets.Item("Garfield Park")
$ws.Range("H6").Value = 26   # was 25
$ws.Range("H7").Value = 35   # was 34

# --- Grand Crossing ---
$ws = $wb.Worksheets.Item("Grand Crossing")
$ws.Range("C6").Value = 21   # was 20
$ws.Range("C7").Value = 25   # was 24

# --- Armour Square ---
$ws = $wb.Worksheets.Item("Armour Square")
$ws.Range("B3").Value = 2   # was 1

# --- Humboldt Park ---
$ws = $wb.Worksheets.Item("Humboldt Park")
$ws.Range("B6").Value = 5   # was 4

# --- Uptown ---
$ws = $wb.Worksheets.Item("Uptown")
$ws.Range("G5").Value = 12   # was 11
$ws.Range("G6").Value = 18   # was 17

# --- Englewood ---
$ws = $wb.Worksheets.Item("Englewood")
$ws.Range("F2").Value = 5   # was 4
$ws.Range("C6").Value = 29   # was 28
$ws.Range("K6").Value = 23   # was 22
$ws.Range("C7").Value = 33   # was 32
$ws.Range("F7").Value = 38   # was 37
$ws.Range("K7").Value = 39   # was 38

# --- By Neighborhood ---
$ws = $wb.Worksheets.Item("By Neighborhood")
$ws.Range("G8").Value = 22   # was 21
$ws.Range("H8").Value = 27   # was 26
$ws.Range("G20").Value = 8   # was 7
$ws.Range("B21").Value = 4   # was 3
$ws.Range("E21").Value = 7   # was 6
$ws.Range("H27").Value = 11   # was 10
$ws.Range("C28").Value = 33   # was 32
$ws.Range("F28").Value = 38   # was 37
$ws.Range("K28").Value = 39   # was 38
$ws.Range("H32").Value = 35   # was 34
$ws.Range("C36").Value = 25   # was 24
$ws.Range("B41").Value = 5   # was 4
$ws.Range("F53").Value = 55   # was 54
$ws.Range("G53").Value = 58   # was 57
$ws.Range("J53").Value = 85   # was 83
$ws.Range("J76").Value = 10   # was 9
$ws.Range("G86").Value = 18   # was 17
$ws.Range("B98").Value = 329   # was 327
$ws.Range("C98").Value = 413   # was 411
$ws.Range("E98").Value = 420   # was 419
$ws.Range("F98").Value = 511   # was 509
$ws.Range("G98").Value = 459   # was 455
$ws.Range("H98").Value = 441   # was 438
$ws.Range("J98").Value = 478   # was 475
$ws.Range("K98").Value = 588   # was 587

# --- Loop ---
$ws = $wb.Worksheets.Item("Loop")
$ws.Range("G3").Value = 13   # was 12
$ws.Range("J3").Value = 23   # was 21
$ws.Range("F6").Value = 41   # was 40
$ws.Range("F7").Value = 55   # was 54
$ws.Range("G7").Value = 58   # was 57
$ws.Range("J7").Value = 85   # was 83

# --- Rogers Park ---
$ws = $wb.Worksheets.Item("Rogers Park")
$ws.Range("J3").Value = 2   # was 1
$ws.Range("J6").Value = 10   # was 9

# --- Edgewater ---
$ws = $wb.Worksheets.Item("Edgewater")
$ws.Range("H4").Value = 11   # was 10
$ws.Range("H5").Value = 11   # was 10

# --- Austin ---
$ws = $wb.Worksheets.Item("Austin")
$ws.Range("G5").Value = 17   # was 16
$ws.Range("H5").Value = 21   # was 20
$ws.Range("G6").Value = 22   # was 21
$ws.Range("H6").Value = 27   # was 26

# --- Chinatown ---
$ws = $wb.Worksheets.Item("Chinatown")
$ws.Range("B6").Value = 3   # was 2
$ws.Range("E6").Value = 3   # was 2
$ws.Range("B7").Value = 4   # was 3
$ws.Range("E7").Value = 7   # was 6

# --- Chicago Lawn ---
$ws = $wb.Worksheets.Item("Chicago Lawn")
$ws.Range("G6").Value = 8   # was 7
$ws.Range("G3").Value = 1   # new cell
